# No semantic content change: this revision only affects how the OOXML
# parts are serialized (attribute ordering) when the authoring tool
# re-saved the template after wiring up the M2Doc-version custom property
# elsewhere in the codebase. The document's text, structure, styles and
# section/page-setup values are identical before and after, so there is
# nothing in the Word object model that needs to change here.
$d = $word.ActiveDocument
